$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 4.75
$ws.Range("O2").Value = 1.73
$ws.Range("P2").Value = 2
$ws.Range("S2").Value = 1.75
$ws.Range("T2").Value = 2.05
$ws.Range("U2").Value = 3
$ws.Range("V2").Value = 1.36
$ws.Range("AA2").Value = 23
$ws.Range("AC2").Value = 4.75
$ws.Range("AZ2").Value = 351
